$wb = $excel.ActiveWorkbook

# Sheet "展览": update F2, F3, F5, F6
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 9333
$ws1.Range("F3").Value = 208
$ws1.Range("F5").Value = 508
$ws1.Range("F6").Value = 460

# Sheet "全部类型": update F2, F3, F5, F7
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9333
$ws4.Range("F3").Value = 208
$ws4.Range("F5").Value = 508
$ws4.Range("F7").Value = 460
